$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 250
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H113").Value = 29084
$ws.Range("I113").Value = 29084
$ws.Range("K113").Value = 29084
$ws.Range("M113").Value = -25830
$ws.Range("H137").Value = 15589.81
$ws.Range("J137").Value = 27000.6
$ws.Range("L137").Value = 81001.79999999999
$ws.Range("N137").Value = -86101.79999999999
$ws.Range("H138").Value = 3474.2454
$ws.Range("I138").Value = 3481.7693
$ws.Range("J138").Value = 3471.8
$ws.Range("K138").Value = 10445.3079
$ws.Range("L138").Value = 10415.4
$ws.Range("M138").Value = -5305.3079
$ws.Range("N138").Value = -20695.4
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7473.7207
$ws.Range("I32").Value = 4519.533
$ws.Range("K32").Value = 4519.533
$ws.Range("M32").Value = -4232.533
$ws.Range("H45").Value = 2112.9412
$ws.Range("I45").Value = 1938.9231
$ws.Range("J45").Value = 2678.5
$ws.Range("K45").Value = 1938.9231
$ws.Range("L45").Value = 2678.5
$ws.Range("M45").Value = -1561.9231
$ws.Range("N45").Value = -3432.5
$ws.Range("H122").Value = 1430777.6
$ws.Range("I122").Value = 1668407.4
$ws.Range("K122").Value = 5005222.199999999
$ws.Range("M122").Value = -5002772.199999999
$ws.Range("H132").Value = 2653646.2
$ws.Range("I132").Value = 5014.55
$ws.Range("J132").Value = 13248173
$ws.Range("K132").Value = 15043.65
$ws.Range("L132").Value = 39744519
$ws.Range("M132").Value = -12513.65
$ws.Range("N132").Value = -39749579
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 103
$ws.Range("I7").Value = 103
$ws.Range("K7").Value = 103
$ws.Range("M7").Value = 10
$ws.Range("H81").Value = 42407.5
$ws.Range("J81").Value = 47890
$ws.Range("L81").Value = 47890
$ws.Range("N81").Value = -50012
$ws.Range("H82").Value = 19200.2
$ws.Range("J82").Value = 45000
$ws.Range("L82").Value = 45000
$ws.Range("N82").Value = -45766
$ws.Range("H84").Value = 42407.5
$ws.Range("J84").Value = 47890
$ws.Range("L84").Value = 143670
$ws.Range("N84").Value = -154278
$ws.Range("H85").Value = 19200.2
$ws.Range("J85").Value = 45000
$ws.Range("L85").Value = 45000
$ws.Range("N85").Value = -47652
$ws.Range("H86").Value = 3890.3684
$ws.Range("I86").Value = 2617.25
$ws.Range("K86").Value = 2617.25
$ws.Range("M86").Value = -1494.25
$ws.Range("H89").Value = 3890.3684
$ws.Range("I89").Value = 2617.25
$ws.Range("K89").Value = 13086.25
$ws.Range("M89").Value = -7470.25
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("H134").Value = 11316.5
$ws.Range("I134").Value = 6752.5884
$ws.Range("J134").Value = 37178.668
$ws.Range("K134").Value = 20257.7652
$ws.Range("L134").Value = 111536.004
$ws.Range("M134").Value = -17722.7652
$ws.Range("N134").Value = -116606.004
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 1031
$ws.Range("I12").Value = 1031
$ws.Range("K12").Value = 1031
$ws.Range("M12").Value = -861
$ws.Range("H31").Value = 64717.258
$ws.Range("I31").Value = 136743.94
$ws.Range("J31").Value = 19700.584
$ws.Range("K31").Value = 136743.94
$ws.Range("L31").Value = 19700.584
$ws.Range("M31").Value = -136448.94
$ws.Range("N31").Value = -20290.584
$ws.Range("H34").Value = 64717.258
$ws.Range("I34").Value = 136743.94
$ws.Range("J34").Value = 19700.584
$ws.Range("K34").Value = 136743.94
$ws.Range("L34").Value = 19700.584
$ws.Range("M34").Value = -136541.94
$ws.Range("N34").Value = -20104.584
$ws.Range("H58").Value = 15137.678
$ws.Range("I58").Value = 5718.6875
$ws.Range("J58").Value = 25184.6
$ws.Range("K58").Value = 5718.6875
$ws.Range("L58").Value = 25184.6
$ws.Range("M58").Value = -5515.6875
$ws.Range("N58").Value = -25590.6
$ws.Range("H86").Value = 8081.476
$ws.Range("I86").Value = 8436.429
$ws.Range("J86").Value = 7371.5713
$ws.Range("K86").Value = 8436.429
$ws.Range("L86").Value = 7371.5713
$ws.Range("M86").Value = -7313.429
$ws.Range("N86").Value = -9617.5713
$ws.Range("H89").Value = 8081.476
$ws.Range("I89").Value = 8436.429
$ws.Range("J89").Value = 7371.5713
$ws.Range("K89").Value = 42182.145
$ws.Range("L89").Value = 36857.85649999999
$ws.Range("M89").Value = -36566.145
$ws.Range("N89").Value = -48089.85649999999
$ws.Range("H134").Value = 43486840
$ws.Range("I134").Value = 2582.3
$ws.Range("J134").Value = 76936270
$ws.Range("K134").Value = 7746.900000000001
$ws.Range("L134").Value = 230808810
$ws.Range("M134").Value = -5211.900000000001
$ws.Range("N134").Value = -230813880
$ws.Range("H136").Value = 15137.678
$ws.Range("I136").Value = 5718.6875
$ws.Range("J136").Value = 25184.6
$ws.Range("K136").Value = 17156.0625
$ws.Range("L136").Value = 75553.79999999999
$ws.Range("M136").Value = -14606.0625
$ws.Range("N136").Value = -80653.79999999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3556.5557
$ws.Range("J68").Value = 3636.0386
$ws.Range("L68").Value = 10908.1158
$ws.Range("N68").Value = -12530.1158
$ws.Range("H71").Value = 3556.5557
$ws.Range("J71").Value = 3636.0386
$ws.Range("L71").Value = 32724.3474
$ws.Range("N71").Value = -40836.3474
$ws.Range("H92").Value = 5159.1113
$ws.Range("I92").Value = 813.8333
$ws.Range("J92").Value = 13849.667
$ws.Range("K92").Value = 2441.4999
$ws.Range("L92").Value = 41549.001
$ws.Range("M92").Value = -1193.4999
$ws.Range("N92").Value = -44045.001
$ws.Range("H107").Value = 4452.75
$ws.Range("I107").Value = 556.26666
$ws.Range("K107").Value = 1668.79998
$ws.Range("M107").Value = 251.20002
$ws.Range("H112").Value = 10608.0625
$ws.Range("I112").Value = 2649.5
$ws.Range("J112").Value = 11745
$ws.Range("K112").Value = 7948.5
$ws.Range("L112").Value = 35235
$ws.Range("M112").Value = -6840.5
$ws.Range("N112").Value = -37451
$ws.Range("H113").Value = 11468.637
$ws.Range("I113").Value = 20380.166
$ws.Range("J113").Value = 774.8
$ws.Range("K113").Value = 61140.49800000001
$ws.Range("L113").Value = 2324.4
$ws.Range("M113").Value = -58970.49800000001
$ws.Range("N113").Value = -6664.4
$ws.Range("H131").Value = 1404.66
$ws.Range("I131").Value = 793.63635
$ws.Range("J131").Value = 1480.1798
$ws.Range("K131").Value = 2380.90905
$ws.Range("L131").Value = 4440.5394
$ws.Range("M131").Value = 2659.09095
$ws.Range("N131").Value = -14520.5394
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7595
$ws.Range("J80").Value = 11491.167
$ws.Range("L80").Value = 11491.167
$ws.Range("N80").Value = -13487.167
$ws.Range("H83").Value = 7595
$ws.Range("J83").Value = 11491.167
$ws.Range("L83").Value = 57455.835
$ws.Range("N83").Value = -67439.83499999999
$ws.Range("H132").Value = 441951.2
$ws.Range("I132").Value = 5401.7085
$ws.Range("K132").Value = 16205.1255
$ws.Range("M132").Value = -13675.1255
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2140.5833
$ws.Range("I16").Value = 2140.5833
$ws.Range("K16").Value = 2140.5833
$ws.Range("M16").Value = -1970.5833
$ws.Range("H40").Value = 8523.409
$ws.Range("I40").Value = 5609.9165
$ws.Range("K40").Value = 5609.9165
$ws.Range("M40").Value = -5473.9165
$ws.Range("H68").Value = 15804.267
$ws.Range("I68").Value = 19028.715
$ws.Range("K68").Value = 19028.715
$ws.Range("M68").Value = -18279.715
$ws.Range("H71").Value = 15804.267
$ws.Range("I71").Value = 19028.715
$ws.Range("K71").Value = 95143.575
$ws.Range("M71").Value = -91399.575
$ws.Range("H110").Value = 34960.43
$ws.Range("J110").Value = 34960.43
$ws.Range("L110").Value = 34960.43
$ws.Range("N110").Value = -43140.43
$ws.Range("H136").Value = 1004992.4
$ws.Range("I136").Value = 26496.334
$ws.Range("K136").Value = 79489.00199999999
$ws.Range("M136").Value = -76939.00199999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 374.5
$ws.Range("I100").Value = 350
$ws.Range("J100").Value = 399
$ws.Range("K100").Value = 700
$ws.Range("L100").Value = 798
$ws.Range("M100").Value = -159
$ws.Range("N100").Value = -1880
$ws.Range("H132").Value = 1784825.2
$ws.Range("I132").Value = 9650.546
$ws.Range("J132").Value = 5335175
$ws.Range("K132").Value = 28951.638
$ws.Range("L132").Value = 16005525
$ws.Range("M132").Value = -26421.638
$ws.Range("N132").Value = -16010585
$ws.Range("H136").Value = 294131.2
$ws.Range("I136").Value = 2713.8667
$ws.Range("J136").Value = 502286.44
$ws.Range("K136").Value = 8141.6001
$ws.Range("L136").Value = 1506859.32
$ws.Range("M136").Value = -5591.6001
$ws.Range("N136").Value = -1511959.32
